# Updates cryptos list price (D) and volume-change (E) columns
# to match the refreshed data snapshot, per commit "Updated cryptos
# list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.993.19"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.05"
$ws.Range("D3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9975"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6253"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9987"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07580"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2917"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.49"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07711"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.839.83"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.947"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6626"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001024"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +19.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.033"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.008.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "226.19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9991"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.181"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9986"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.476"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1375"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.487"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.087"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.003"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.185"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05230"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.838"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7333"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.684"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.234.36"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.753"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01778"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.320"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8976"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9987"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.76"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.977.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5099"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4029"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.879"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05735"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.657"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.21%  "
